# Refresh the crypto price/volume data cells on the active sheet to match
# the latest scrape, keeping cell types as plain text (matching source file).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.966.93'
$ws.Range("E2").Value = '  +1.08%  '
$ws.Range("D3").Value = '3.520.80'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.91%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.84'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.78%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.520.53'
$ws.Range("E8").Value = '  +0.23%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.596'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.140'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.89%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.17'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.68%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.439'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.45%  '
$ws.Range("D13").Value = '4.130.22'
$ws.Range("E13").Value = '  +0.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.23'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +10.04%  '
$ws.Range("E15").Value = '  +1.07%  '
$ws.Range("D16").Value = '67.940.34'
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("D18").Value = '3.531.53'
$ws.Range("E18").Value = '  +0.01%  '
$ws.Range("E19").Value = '  +0.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '401.34'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '74.07'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.27%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.545'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.19%  '
$ws.Range("E26").Value = '  +0.77%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000123'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.48%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.53'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.32%  '
$ws.Range("E29").Value = '  -2.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("E31").Value = '  -0.29%  '
$ws.Range("E32").Value = '  -1.00%  '
$ws.Range("E33").Value = '  +1.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.98'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.52'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.51%  '
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("E37").Value = '  -2.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '163.34'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.883'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.81%  '
$ws.Range("E40").Value = '  +0.85%  '
$ws.Range("E41").Value = '  +8.04%  '
$ws.Range("E42").Value = '  -1.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.72'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").Value = '2.892.13'
$ws.Range("E44").Value = '  +1.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '26.51'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.46%  '
$ws.Range("E46").Value = '  -2.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.95'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.00%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '42.51'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '351.66'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.51%  '
$ws.Range("E50").Value = '  +0.33%  '
$ws.Range("E51").Value = '  -1.05%  '
